$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Column A (PackageName) for all data rows to new package name
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "ExePkgSep30T2"
}

# Column K (childName) values are replaced with the childVariableName (column L) value,
# effectively removing the distinct "childName" content for these rows.
$rowsToSync = 3,4,5,6,7,8,12,13,14,15
foreach ($r in $rowsToSync) {
    $ws.Cells.Item($r, 11).Value2 = $ws.Cells.Item($r, 12).Value2
}
